$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Reference sheet used only to borrow the bold/bordered header style
# (style index 2 in styles.xml) and the matching blank-index style
# (style index 1 applied via "A2") that every other data sheet uses.
# ------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item("trafo3w_std_types")

# ------------------------------------------------------------------
# New sheet: pwl_cost (appended after trafo3w_std_types)
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$pwlCost = $wb.Worksheets.Add($null, $lastSheet)
$pwlCost.Name = "pwl_cost"

$pwlCost.Range("B1").Value = "power_type"
$pwlCost.Range("C1").Value = "element"
$pwlCost.Range("D1").Value = "et"
$pwlCost.Range("E1").Value = "points"

$styleSrc.Range("B1").Copy()
$pwlCost.Range("B1:E1").PasteSpecial(-4122)

$styleSrc.Range("A2").Copy()
$pwlCost.Range("A2").PasteSpecial(-4122)

$pwlCost.Range("Q20").Select()

# ------------------------------------------------------------------
# New sheet: poly_cost (appended after pwl_cost)
# ------------------------------------------------------------------
$polyCost = $wb.Worksheets.Add($null, $pwlCost)
$polyCost.Name = "poly_cost"

$polyCost.Range("B1").Value = "element"
$polyCost.Range("C1").Value = "et"
$polyCost.Range("D1").Value = "cp0_eur"
$polyCost.Range("E1").Value = "cp1_eur_per_mw"
$polyCost.Range("F1").Value = "cp2_eur_per_mw2"
$polyCost.Range("G1").Value = "cq0_eur"
$polyCost.Range("H1").Value = "cq1_eur_per_mvar"
$polyCost.Range("I1").Value = "cq2_eur_per_mvar2"

$styleSrc.Range("B1").Copy()
$polyCost.Range("B1:I1").PasteSpecial(-4122)

$styleSrc.Range("A2").Copy()
$polyCost.Range("A2").PasteSpecial(-4122)

$polyCost.Range("Q20").Select()

# ------------------------------------------------------------------
# "general" sheet: add a use_opf column (new opf analysis flag)
# ------------------------------------------------------------------
$general = $wb.Worksheets.Item("general")
$general.Range("C1").Value = "use_opf"
$general.Range("C2").Value = $false

# ------------------------------------------------------------------
# Restore the previously-selected cell on trafo3w_std_types (it is no
# longer the active tab) and make "general" the active tab instead.
# ------------------------------------------------------------------
$styleSrc.Range("K27").Select()

$general.Activate()
$general.Range("E15").Select()
